# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "297.16";       E = "-1.58%" }
    3  = @{ D = "31.35";        E = "-0.37%" }
    4  = @{ D = "5.088";        E = "-1.10%" }
    5  = @{ D = "0.08022";      E = "8.56%" }
    6  = @{ D = "2.639";        E = "43.98%" }
    7  = @{ E = "-0.77%" }
    8  = @{ E = "1.02%" }
    9  = @{ D = "0.9284";       E = "-0.05%" }
    10 = @{ D = "0.1741";       E = "1.95%" }
    11 = @{ D = "0.07194";      E = "-1.37%" }
    12 = @{ D = "0.08980";      E = "10.99%" }
    13 = @{ D = "0.03083";      E = "1.48%" }
    14 = @{ D = "0.1001";       E = "0.69%" }
    15 = @{ D = "0.001503";     E = "0.11%" }
    16 = @{ D = "0.006018";     E = "-0.46%" }
    17 = @{ D = "3.545";        E = "2.14%" }
    18 = @{ E = "0.97%" }
    19 = @{ E = "-0.88%" }
    20 = @{ D = "0.1342";       E = "1.96%" }
    21 = @{ D = "3.920";        E = "-14.39%" }
    22 = @{ D = "0.1620";       E = "2.57%" }
    23 = @{ D = "0.04581";      E = "-1.46%" }
    24 = @{ E = "1.94%" }
    25 = @{ D = "0.004420";     E = "-7.57%" }
    26 = @{ D = "0.0001199";    E = "-7.53%" }
    27 = @{ D = "0.0003417";    E = "82.55%" }
    39 = @{ D = "0.01773";      E = "2.64%" }
    40 = @{ D = "0.04479";      E = "-0.67%" }
    41 = @{ D = "0.006869";     E = "-3.03%" }
    42 = @{ D = "0.1344";       E = "-0.41%" }
    43 = @{ D = "0.002134";     E = "-0.47%" }
    44 = @{ D = "0.009836";     E = "-6.21%" }
    45 = @{ D = "0.00006525";   E = "5.07%" }
    46 = @{ D = "0.00000000749"; E = "-0.14%" }
    47 = @{ D = "0.008723";     E = "-14.60%" }
    48 = @{ E = "-55.51%" }
    49 = @{ D = "0.00002097";   E = "-0.14%" }
    50 = @{ D = "0.0001997";    E = "-0.07%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
